$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LeetCode problems to append: (S.No, Name, C-style-source-row, D-style-source-row)
# Row 40 -> Problem 101 Symmetric Tree            (C/D like row 2:  C=Good/border7, D=Good/border8)
# Row 41 -> Problem 104 Maximum Depth of Binary Tree   (same as row 2)
# Row 42 -> Problem 107 Binary Tree Level Order Traversal II (C/D like row 9: C=Good, D=Neutral)
# Row 43 -> Problem 108 Convert Sorted Array to Binary Search Tree (same as row 2)
# Row 44 -> Problem 110 Balanced Binary Tree       (same as row 9)

$rows = @(
    @{ Row = 40; Num = 101; Name = "Symmetric Tree";                               StyleRow = 2 },
    @{ Row = 41; Num = 104; Name = "Maximum Depth of Binary Tree";                 StyleRow = 2 },
    @{ Row = 42; Num = 107; Name = "Binary Tree Level Order Traversal II";         StyleRow = 9 },
    @{ Row = 43; Num = 108; Name = "Convert Sorted Array to Binary Search Tree";   StyleRow = 2 },
    @{ Row = 44; Num = 110; Name = "Balanced Binary Tree";                         StyleRow = 9 }
)

foreach ($r in $rows) {
    $destRow = $r.Row

    # Copy A/B formatting from the last existing data row (row 39)
    $ws.Range("A39:B39").Copy()
    $ws.Range("A" + $destRow + ":B" + $destRow).PasteSpecial(-4122)

    # Copy C/D formatting from a row with the matching Good/Neutral style combo
    $srcRow = $r.StyleRow
    $ws.Range("C" + $srcRow + ":D" + $srcRow).Copy()
    $ws.Range("C" + $destRow + ":D" + $destRow).PasteSpecial(-4122)

    # Set the actual values
    $ws.Range("A" + $destRow).Value = $r.Num
    $ws.Range("B" + $destRow).Value = $r.Name
}

$excel.CutCopyMode = 0

# Update selection to match the new active cell recorded in the saved file
$ws.Range("K36").Select()
